$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Aggiornamento_0") gets refreshed workflow/document identifiers.
# Values are written left-to-right so newly introduced shared strings are
# registered in the same order as the source data extract.
$ws.Range("A2").Value = "Aggiornamento_0"
$ws.Range("B2").Value = "REGIONE_CAMPANIA"
$ws.Range("C2").Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Range("D2").Value = "765080b208dfddbfff8fc6512f911b0817499f80292a6751c453920224643978.812baa6b68^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721299538721"
$ws.Range("F2").Value = "18-07-2024:12:46:35"

# Row 3 now holds the "Creazione_0" event (previously the 4th row), with
# refreshed values; the old "Creazione_1" row is gone.
$ws.Range("A3").Value = "Creazione_0"
$ws.Range("B3").Value = "REGIONE_CAMPANIA"
$ws.Range("C3").Value = "NGNVCN92S19L259C^^^&2.16.840.1.113883.2.9.4.3.2&ISO"
$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.05e2ca76c0^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721299538721"
$ws.Range("F3").Value = "18-07-2024:12:45:40"

# Remove the old row 4 (former "Creazione_0" row), leaving only 3 data rows
$ws.Range("A4:F4").Delete()
